$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'258.25"
$ws.Range("E2").Value = "'5.35%"
$ws.Range("D3").Value = "'28.01"
$ws.Range("E3").Value = "'-4.01%"
$ws.Range("D4").Value = "'5.208"
$ws.Range("E4").Value = "'-0.96%"
$ws.Range("D5").Value = "'0.05907"
$ws.Range("E5").Value = "'3.52%"
$ws.Range("D6").Value = "'6.706"
$ws.Range("E6").Value = "'1.36%"
$ws.Range("D7").Value = "'0.8733"
$ws.Range("E7").Value = "'2.69%"
$ws.Range("D8").Value = "'0.9788"
$ws.Range("E8").Value = "'13.87%"
$ws.Range("D9").Value = "'0.1413"
$ws.Range("E9").Value = "'3.06%"
$ws.Range("D10").Value = "'0.07187"
$ws.Range("E10").Value = "'2.21%"
$ws.Range("D11").Value = "'0.03162"
$ws.Range("E11").Value = "'-1.28%"
$ws.Range("D12").Value = "'0.09237"
$ws.Range("E12").Value = "'-0.50%"
$ws.Range("D13").Value = "'0.001544"
$ws.Range("E13").Value = "'1.48%"
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").Value = "'0.006019"
$ws.Range("E14").Value = "'-0.54%"
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D15").Value = "'3.494"
$ws.Range("E15").Value = "'0.04%"
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D16").Value = "'3.231"
$ws.Range("E16").Value = "'1.24%"
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").Value = "'2.204"
$ws.Range("E17").Value = "'-1.14%"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.01062"
$ws.Range("E18").Value = "'5.68%"
$ws.Range("E19").Value = "'-1.18%"
$ws.Range("D20").Value = "'0.03722"
$ws.Range("E20").Value = "'12.84%"
$ws.Range("D21").Value = "'0.1298"
$ws.Range("E21").Value = "'1.68%"
$ws.Range("D22").Value = "'3.530"
$ws.Range("E22").Value = "'0.85%"
$ws.Range("D23").Value = "'0.04202"
$ws.Range("E23").Value = "'2.68%"
$ws.Range("D24").Value = "'0.1363"
$ws.Range("E24").Value = "'-1.21%"
$ws.Range("D25").Value = "'0.001219"
$ws.Range("E25").Value = "'-0.15%"
$ws.Range("D26").Value = "'0.004545"
$ws.Range("E26").Value = "'9.74%"
$ws.Range("D27").Value = "'0.0001199"
$ws.Range("E27").Value = "'-0.10%"
$ws.Range("E28").Value = "'33.61%"
$ws.Range("D40").Value = "'0.03835"
$ws.Range("E40").Value = "'2.07%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.005478"
$ws.Range("E41").Value = "'6.23%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1103"
$ws.Range("E42").Value = "'3.71%"
$ws.Range("D43").Value = "'0.002299"
$ws.Range("E43").Value = "'-6.15%"
$ws.Range("E44").Value = "'13.73%"
$ws.Range("D45").Value = "'0.00005418"
$ws.Range("E45").Value = "'2.96%"
$ws.Range("E46").Value = "'-0.10%"
$ws.Range("D47").Value = "'0.08546"
$ws.Range("E47").Value = "'13.88%"
$ws.Range("E48").Value = "'-12.53%"
$ws.Range("E49").Value = "'-0.10%"
$ws.Range("E50").Value = "'-0.10%"
